$d = $word.ActiveDocument

# --- Change 1: "a pretty decent amount" -> "a decent amount" ---
$null = $d.Content.Find.Execute(
    "a pretty decent amount of online documentation.", $true, $false, $false, $false, $false,
    $true, 1, $false, "a decent amount of online documentation.", 2)

# --- Change 2: drop the redundant "also" before "relatively easy" ---
$null = $d.Content.Find.Execute(
    "current user was also relatively easy.", $true, $false, $false, $false, $false,
    $true, 1, $false, "current user was relatively easy.", 2)

# --- Change 3: expand the Authlogic "simple if statement" sentence with new material ---
$oldTail = "a simple if statement (very similar to Devise)."
$newTail = "a simple if statement. Authlogic did require you to create more code yourself (there is a user_sessions controller that we defined, although this does give you a better understanding of what the code is actually doing). All of the user creation and authentication was simple in Authlogic as well (Rails in general seems to do a very good job at simplifying stuff)."
$null = $d.Content.Find.Execute(
    $oldTail, $true, $false, $false, $false, $false,
    $true, 1, $false, $newTail, 2)

# --- Split that paragraph: give the _GoBack bookmark its own empty paragraph ---
$splitPoint = $d.Content
$null = $splitPoint.Find.Execute("very good job at simplifying stuff).")
$splitPoint.Collapse(0)
$splitPoint.InsertBefore([char]13)

# --- Remove the two now-orphaned blank paragraphs that used to follow the bookmark ---
$bm = $d.Bookmarks.Item("_GoBack")
$bmParaIndex = $bm.Range.Paragraphs.Item(1).Index
$d.Paragraphs.Item($bmParaIndex + 1).Range.Delete()
$d.Paragraphs.Item($bmParaIndex + 1).Range.Delete()

# --- Change 4: replace the "We chose to select…" placeholder with the real summary ---
$finalText = "Both of these packages seemed to work very well when creating an app and were pretty much equally as effective at implementing authentication, but after discussing and weighing our options, we decided to go with Authlogic. Both Authlogic and Devise seemed like a good choice (and Devise may have won if we were starting an app from scratch since having admin privileges seemed easy to do) but Authlogic seemed to be easier to implement and we believe it will be easier to incorporate into an already existing project. "
$null = $d.Content.Find.Execute(
    "We chose to select" + [char]0x2026, $true, $false, $false, $false, $false,
    $true, 1, $false, $finalText, 2)
